$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week's price record was added to the historical table. Excel's
# native row-insert shifts row 31 and everything below it down by one
# row (carrying values/styles along), matching the diff where old row
# r (31..95) becomes new row r+1 (32..96).
$ws.Rows.Item(31).Insert()

# Populate the newly-inserted row 31 with this week's record.
$ws.Range("A31").Value = 1
$ws.Range("B31").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C31").Value = "Arica y Parinacota"
$ws.Range("D31").Value = 45177
$ws.Range("E31").Value = 15
$ws.Range("F31").Value = 100112031
$ws.Range("G31").Value = "Poroto verde"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 1300
$ws.Range("K31").Value = 1000
$ws.Range("L31").Value = 1200
$ws.Range("M31").Value = 1100
$ws.Range("N31").Value = '$/kilo'
$ws.Range("O31").Value = "Región de Arica y Parinacota"
$ws.Range("P31").Value = 1100
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = "Hortaliza"
